# WalletLedger update — data updated till 17 Jan 2021 8AM
# Row 20's debited amount was corrected, and two new ledger rows
# (21: a manually-added credit, 22: a new ordered-amount debit) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: correct the debited amount, and drop the italic "last row"
#     emphasis from D20 since row 20 is no longer the final entry ---
$ws.Range("B20").Value = 46800
$ws.Range("D20").Font.Italic = $false

# --- Row 21 (new): 16-Jan credit, manually added balance ---
$ws.Range("A20").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 43846
$ws.Range("C16").Copy($ws.Range("C21"))
$ws.Range("C21").Value = 369564
$ws.Range("D16").Copy($ws.Range("D21"))
$ws.Range("D21").Value = "Manual Added"
$ws.Range("E21").Formula = '=IF(A21="","",SUM(E20-B21+C21))'

# --- Row 22 (new): 17-Jan debit, ordered amount — now the final row, so it
#     takes over the italic emphasis style previously on D20 ---
$ws.Range("A20").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 43847
$ws.Range("B19").Copy($ws.Range("B22"))
$ws.Range("B22").Value = 13520
$ws.Range("D19").Copy($ws.Range("D22"))
$ws.Range("D22").Value = "Ordered Amount"
$ws.Range("D22").Font.Italic = $true
$ws.Range("E22").Formula = '=IF(A22="","",SUM(E21-B22+C22))'

# --- Move the active selection to the new last entry ---
$ws.Range("D22").Select() | Out-Null
